$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.781.36'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.629.14'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.77%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.39'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.66'
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').Value = '  +1.27%  '
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.853.62'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '1.621.30'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.85'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '25.758.12'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.24'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.93'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.28'
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.997'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.92'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('E27').Value = '  +3.07%  '
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.50'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.38'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D37').Value = '1.141.91'
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.543'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.50'
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.996'
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.51'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.59'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.69'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.805'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').Value = '1.764.05'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.30'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.46'
$ws.Range('E48').Value = '  +6.77%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0511'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.33'
$ws.Range('E51').Value = '  -1.81%  '
